$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1.348535280068247
$ws.Range("C3").Value = 2.280027247411239
$ws.Range("D3").Value = 3.618982841204426
$ws.Range("E3").Value = 5.427926478372039
$ws.Range("F3").Value = 7.745174434347612
$ws.Range("B4").Value = 1.002720984225978
$ws.Range("C4").Value = 1.746133973763164
$ws.Range("D4").Value = 2.849455207820629
$ws.Range("E4").Value = 4.383495450832996
$ws.Range("F4").Value = 6.398574505914165
$ws.Range("B5").Value = 0.7363758790605551
$ws.Range("C5").Value = 1.32073792052963
$ws.Range("D5").Value = 2.21689655653529
$ws.Range("E5").Value = 3.500746455980269
$ws.Range("F5").Value = 5.232592204231606
$ws.Range("B6").Value = 0.534354214691467
$ws.Range("C6").Value = 0.9868075587613082
$ws.Range("D6").Value = 1.704160804274016
$ws.Range("E6").Value = 2.764096942459139
$ws.Range("F6").Value = 4.23429477656448
$ws.Range("B7").Value = 0.3833978806917301
$ws.Range("C7").Value = 0.7285524660266689
$ws.Range("D7").Value = 1.294433490935827
$ws.Range("E7").Value = 2.157425595591394
$ws.Range("F7").Value = 3.389634940376918
